# Auto-generated edit script applying numeric cell updates to the Titan_Profits workbook
# (ALC/ARM/BSM/CRP/GSM/LTW/WVR sheets) per the target diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1235021.2
$ws.Range("I28").Value = 1389148.9
$ws.Range("J28").Value = 2000
$ws.Range("K28").Value = 1389148.9
$ws.Range("L28").Value = 2000
$ws.Range("M28").Value = -1388663.9
$ws.Range("N28").Value = -2970
$ws.Range("H127").Value = 792.7143
$ws.Range("I127").Value = 225
$ws.Range("J127").Value = 1019.8
$ws.Range("K127").Value = 675
$ws.Range("L127").Value = 3059.4
$ws.Range("M127").Value = 4285
$ws.Range("N127").Value = -12979.4
$ws.Range("H137").Value = 41668090
$ws.Range("I137").Value = 62500988
$ws.Range("J137").Value = 2284.375
$ws.Range("K137").Value = 187502964
$ws.Range("L137").Value = 6853.125
$ws.Range("M137").Value = -187500414
$ws.Range("N137").Value = -11953.125
$ws.Range("H138").Value = 4210913.5
$ws.Range("I138").Value = 1110651.5
$ws.Range("K138").Value = 3331954.5
$ws.Range("M138").Value = -3326814.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19582.574
$ws.Range("I32").Value = 3025.8135
$ws.Range("J32").Value = 508007
$ws.Range("K32").Value = 3025.8135
$ws.Range("L32").Value = 508007
$ws.Range("M32").Value = -2738.8135
$ws.Range("N32").Value = -508581
$ws.Range("H74").Value = 3334.0166
$ws.Range("I74").Value = 1062.4694
$ws.Range("J74").Value = 13452.728
$ws.Range("K74").Value = 1062.4694
$ws.Range("L74").Value = 13452.728
$ws.Range("M74").Value = -188.4694
$ws.Range("N74").Value = -15200.728
$ws.Range("H77").Value = 3334.0166
$ws.Range("I77").Value = 1062.4694
$ws.Range("J77").Value = 13452.728
$ws.Range("K77").Value = 5312.347
$ws.Range("L77").Value = 67263.64
$ws.Range("M77").Value = -944.3469999999998
$ws.Range("N77").Value = -75999.64

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H139").Value = 74530
$ws.Range("J139").Value = 74530
$ws.Range("L139").Value = 74530
$ws.Range("N139").Value = -84810

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2238.0256
$ws.Range("I31").Value = 1344.7084
$ws.Range("K31").Value = 1344.7084
$ws.Range("M31").Value = -1049.7084
$ws.Range("H34").Value = 2238.0256
$ws.Range("I34").Value = 1344.7084
$ws.Range("K34").Value = 1344.7084
$ws.Range("M34").Value = -1142.7084
$ws.Range("H58").Value = 1359.902
$ws.Range("I58").Value = 866
$ws.Range("J58").Value = 2803.6155
$ws.Range("K58").Value = 866
$ws.Range("L58").Value = 2803.6155
$ws.Range("M58").Value = -663
$ws.Range("N58").Value = -3209.6155
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()  # was -7560
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()  # was -21306
$ws.Range("H68").Value = 27333.334
$ws.Range("J68").Value = 27333.334
$ws.Range("L68").Value = 27333.334
$ws.Range("N68").Value = -28831.334
$ws.Range("H69").Value = 15000
$ws.Range("I69").Value = 15000
$ws.Range("K69").Value = 15000
$ws.Range("M69").Value = -14251
$ws.Range("H70").Value = 29000
$ws.Range("J70").Value = 29000
$ws.Range("L70").Value = 29000
$ws.Range("N70").Value = -29630
$ws.Range("H71").Value = 27333.334
$ws.Range("J71").Value = 27333.334
$ws.Range("L71").Value = 82000.00199999999
$ws.Range("N71").Value = -89488.00199999999
$ws.Range("H72").Value = 15000
$ws.Range("I72").Value = 15000
$ws.Range("K72").Value = 45000
$ws.Range("M72").Value = -41256
$ws.Range("H73").Value = 29000
$ws.Range("J73").Value = 29000
$ws.Range("L73").Value = 29000
$ws.Range("N73").Value = -31184
$ws.Range("H111").Value = 54500.668
$ws.Range("J111").Value = 54500.668
$ws.Range("L111").Value = 54500.668
$ws.Range("N111").Value = -62680.668
$ws.Range("H119").Value = 41104.4
$ws.Range("J119").Value = 41104.4
$ws.Range("L119").Value = 41104.4
$ws.Range("N119").Value = -50780.4
$ws.Range("H136").Value = 1359.902
$ws.Range("I136").Value = 866
$ws.Range("J136").Value = 2803.6155
$ws.Range("K136").Value = 2598
$ws.Range("L136").Value = 8410.8465
$ws.Range("M136").Value = -48
$ws.Range("N136").Value = -13510.8465

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 50000
$ws.Range("J4").Value = 50000
$ws.Range("L4").Value = 50000
$ws.Range("N4").Value = -50224
$ws.Range("H12").Value = 5000
$ws.Range("J12").Value = 5000
$ws.Range("L12").Value = 5000
$ws.Range("N12").Value = -5280
$ws.Range("H123").Value = 9704.125
$ws.Range("J123").Value = 9704.125
$ws.Range("L123").Value = 9704.125
$ws.Range("N123").Value = -14604.125
$ws.Range("H132").Value = 2295.7793
$ws.Range("I132").Value = 2028.8334
$ws.Range("J132").Value = 3325.4285
$ws.Range("K132").Value = 6086.5002
$ws.Range("L132").Value = 9976.2855
$ws.Range("M132").Value = -3556.5002
$ws.Range("N132").Value = -15036.2855

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 3000
$ws.Range("J18").Value = 3000
$ws.Range("L18").Value = 3000
$ws.Range("N18").Value = -3344
$ws.Range("H136").Value = 3053.6167
$ws.Range("I136").Value = 1756
$ws.Range("J136").Value = 5648.85
$ws.Range("K136").Value = 5268
$ws.Range("L136").Value = 16946.55
$ws.Range("M136").Value = -2718
$ws.Range("N136").Value = -22046.55

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 31938.5
$ws.Range("J109").Value = 31938.5
$ws.Range("L109").Value = 31938.5
$ws.Range("N109").Value = -34712.5
$ws.Range("H123").Value = 28130.303
$ws.Range("J123").Value = 28130.303
$ws.Range("L123").Value = 28130.303
$ws.Range("N123").Value = -37930.303
$ws.Range("H132").Value = 8476662
$ws.Range("I132").Value = 13159935
$ws.Range("J132").Value = 2167.5715
$ws.Range("K132").Value = 39479805
$ws.Range("L132").Value = 6502.7145
$ws.Range("M132").Value = -39477275
$ws.Range("N132").Value = -11562.7145
$ws.Range("H133").Value = 69905
$ws.Range("J133").Value = 69905
$ws.Range("L133").Value = 69905
$ws.Range("N133").Value = -80025

Write-Output "Applied 161 cell updates across 7 sheets."